# Fill in the "service_time_mins" (column AA) values that were previously
# blank for every stop, now that vrp_solve_07.py resolves them per-route:
#   - rows 2-16  (ALS stops, column Z = "ALS")  -> 5
#   - rows 17-31 (BLS stops, column Z = "BLS")  -> 7

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("AA2:AA16").Value = 5
$ws.Range("AA17:AA31").Value = 7

# Leave the sheet positioned/selected the way it was when the author left
# off reviewing the newly-populated column.
[void]$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 24
[void]$ws.Range("AA2:AA16").Select()
